$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update degradation curve / DC RTE values
$ws.Range("B1").Value = 0.003
$ws.Range("B3").Value = 0.004

# Update the active selection to B4
$ws.Range("B4").Select()
